$wb = $excel.ActiveWorkbook

# --- Income sheet: row 4, column G (Total) ---
$wsIncome = $wb.Worksheets.Item("Income")
$wsIncome.Range("G4").Value = 300

# --- Expenditure sheet: rows 4-21, column L (Total) ---
$wsExpenditure = $wb.Worksheets.Item("Expenditure")

$wsExpenditure.Range("L4").Value = 10
$wsExpenditure.Range("L5").Value = 5000
$wsExpenditure.Range("L6").Value = 10000
$wsExpenditure.Range("L7").Value = ""
$wsExpenditure.Range("L8").Value = 10
$wsExpenditure.Range("L9").Value = ""
$wsExpenditure.Range("L10").Value = 300
$wsExpenditure.Range("L11").Value = ""
$wsExpenditure.Range("L12").Value = ""
$wsExpenditure.Range("L13").Value = ""
$wsExpenditure.Range("L14").Value = ""
$wsExpenditure.Range("L15").Value = ""
$wsExpenditure.Range("L16").Value = ""
$wsExpenditure.Range("L17").Value = ""
$wsExpenditure.Range("L18").Value = ""
$wsExpenditure.Range("L19").Value = 10
$wsExpenditure.Range("L20").Value = ""
$wsExpenditure.Range("L21").Value = ""
